# "saida do voiter do conglomerado master"
# VOITER leaves the MASTER/WILL/VOITER/LETSBANK conglomerate and becomes its
# own entry ("BANCO PLENO S.A. (VOITER)") with two new alias rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 54: VOITER used to roll up under the MASTER conglomerate group
#     (A/C = "CONGLOMERADO: MASTER, WILL, VOITER E LETSBANK"); now it rolls
#     up under its own "BANCO PLENO S.A. (VOITER)" group. Column B keeps the
#     existing "VOITER" alias untouched. Unhide first so the engine doesn't
#     stamp a spurious auto-fit row height on a write into a hidden row.
$ws.Range("A54").EntireRow.Hidden = $false
$ws.Range("A54").Value = "BANCO PLENO S.A. (VOITER)"
$ws.Range("C54").Value = "BANCO PLENO S.A. (VOITER)"
$ws.Range("A54").Font.Color = 0
$ws.Range("C54").Font.Color = 0

# --- Two new alias rows for the new "BANCO PLENO S.A. (VOITER)" group.
$ws.Range("A197").Value = "BANCO PLENO S.A. (VOITER)"
$ws.Range("B197").Value = "BANCO PLENO"
$ws.Range("C197").Value = "BANCO PLENO S.A. (VOITER)"
$ws.Range("A197").Font.Color = 0
$ws.Range("C197").Font.Color = 0

$ws.Range("A198").Value = "BANCO PLENO S.A. (VOITER)"
$ws.Range("B198").Value = "BANCO PLENO S.A."
$ws.Range("C198").Value = "BANCO PLENO S.A. (VOITER)"
$ws.Range("A198").Font.Color = 0
$ws.Range("C198").Font.Color = 0

# --- Rebuild the AutoFilter over the new data extent (A1:C196) with its
#     single discrete-value filter now pointing at the new conglomerate
#     label. This recomputes every row's hidden state from the filter match
#     (unhiding the remaining MASTER-conglomerate rows 60/62/182, re-hiding
#     the MERCADO CREDITO rows 138/193/194/196 that no longer match, etc.)
$ws.AutoFilterMode = $false
[void]($ws.Range("A1:C196").AutoFilter(1, @("CONGLOMERADO: MASTER, WILL, VOITER E LETSBANK"), 7))

# --- Rows 54/197/198 hold the new "BANCO PLENO S.A. (VOITER)" label, which
#     does not match the active filter criteria, but the author left them
#     visible (manually unhidden), so force them back to visible.
$ws.Range("A54").EntireRow.Hidden = $false
$ws.Range("A197").EntireRow.Hidden = $false
$ws.Range("A198").EntireRow.Hidden = $false

# --- Keep the _FilterDatabase defined name in sync with the new filter range.
foreach ($n in $wb.Names) {
    if ($n.Name() -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$C`$196"
    }
}

# --- Move the selection cursor to match where the author ended up.
$ws.Range("B208").Select()
